$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.288.41"
$ws.Range("E2").Value = "  +4.29%  "
$ws.Range("D3").Value = "3.054.71"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'549.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.57%  "
$ws.Range("D6").Value = "'139.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.83%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.047.99"
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("D10").Value = "'6.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("E13").Value = "  +3.86%  "
$ws.Range("D14").Value = "'34.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.59%  "
$ws.Range("D15").Value = "3.556.85"
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").Value = "63.308.12"
$ws.Range("E16").Value = "  +4.20%  "
$ws.Range("D17").Value = "3.060.80"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").Value = "'6.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").Value = "'480.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.08%  "
$ws.Range("D21").Value = "'13.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("D22").Value = "'0.673"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("D23").Value = "'7.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.07%  "
$ws.Range("D24").Value = "'80.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'12.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.38%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +4.38%  "
$ws.Range("D28").Value = "'7.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +7.58%  "
$ws.Range("E31").Value = "  +2.99%  "
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("E33").Value = "  +7.59%  "
$ws.Range("E34").Value = "  +7.15%  "
$ws.Range("D35").Value = "'55.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("E36").Value = "  +3.85%  "
$ws.Range("D37").Value = "'461.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("E38").Value = "  +4.65%  "
$ws.Range("D39").Value = "3.115.73"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("E40").Value = "  +4.09%  "
$ws.Range("D41").Value = "'0.117"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").Value = "'8.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.55%  "
$ws.Range("D43").Value = "'2.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.61%  "
$ws.Range("D44").Value = "'28.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.82%  "
$ws.Range("D45").Value = "'0.251"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  +4.91%  "
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").Value = "'115.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("E50").Value = "  +2.98%  "
$ws.Range("E51").Value = "  +5.32%  "
